# Update the "Förändrad" date column (C) for rows 2-11 from 45175 (2023-09-06)
# to 45183 (2023-09-14), keeping the existing date formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C11").Value = 45183
